$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MISSISSAUGA"
$ws.Range("B2").Value = "Tough Bud Cannabis"
$ws.Range("C2").Value = "296 LAKESHORE RD W"
$ws.Range("E2").Value = "https://toughbud.ca/shop-missisauga/"
$ws.Range("F2").Value = "Buddi"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "['Curbside pickup', 'Delivery', 'In-store pickup', 'In-store shopping', 'Same-day delivery']"
$ws.Range("I2").Value = 19052782222
$ws.Range("J2").Value = "['Delivery serve within 30 km radius', 'Same-day delivery']"
$ws.Range("K2").Value = '$40 minimum order '
